# Auto-generated Excel COM-interop script to apply the Ifrit_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 891.7727
$ws.Range("I28").Value = 1117
$ws.Range("J28").Value = 497.625
$ws.Range("K28").Value = 1117
$ws.Range("L28").Value = 497.625
$ws.Range("M28").Value = -632
$ws.Range("N28").Value = -1467.625
$ws.Range("H62").Value = 2521.139
$ws.Range("I62").Value = 1782.5
$ws.Range("J62").Value = 2668.8667
$ws.Range("K62").Value = 1782.5
$ws.Range("L62").Value = 2668.8667
$ws.Range("M62").Value = -1158.5
$ws.Range("N62").Value = -3916.8667
$ws.Range("H65").Value = 2521.139
$ws.Range("I65").Value = 1782.5
$ws.Range("J65").Value = 2668.8667
$ws.Range("K65").Value = 8912.5
$ws.Range("L65").Value = 13344.3335
$ws.Range("M65").Value = -5792.5
$ws.Range("N65").Value = -19584.3335
$ws.Range("H88").Value = 3025.3333
$ws.Range("I88").Value = 1200
$ws.Range("J88").Value = 3390.4
$ws.Range("K88").Value = 1200
$ws.Range("L88").Value = 3390.4
$ws.Range("M88").Value = -794
$ws.Range("N88").Value = -4202.4
$ws.Range("H91").Value = 3025.3333
$ws.Range("I91").Value = 1200
$ws.Range("J91").Value = 3390.4
$ws.Range("K91").Value = 1200
$ws.Range("L91").Value = 3390.4
$ws.Range("M91").Value = 204
$ws.Range("N91").Value = -6198.4
$ws.Range("H94").Value = 5390
$ws.Range("I94").Value = 4322.222
$ws.Range("K94").Value = 4322.222
$ws.Range("M94").Value = -3871.222
$ws.Range("H107").Value = 1059.9166
$ws.Range("I107").Value = 1283.4546
$ws.Range("K107").Value = 1283.4546
$ws.Range("M107").Value = 636.5454
$ws.Range("H113").Value = 1970.3636
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 1741.5555
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 1741.5555
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -8249.5555
$ws.Range("H129").Value = 1626.1818
$ws.Range("I129").Value = 1051.375
$ws.Range("J129").Value = 3159
$ws.Range("K129").Value = 3154.125
$ws.Range("L129").Value = 9477
$ws.Range("M129").Value = 1845.875
$ws.Range("N129").Value = -19477

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3648
$ws.Range("I88").Value = 4600
$ws.Range("J88").Value = 2696
$ws.Range("K88").Value = 4600
$ws.Range("L88").Value = 2696
$ws.Range("M88").Value = -4194
$ws.Range("N88").Value = -3508
$ws.Range("H91").Value = 3648
$ws.Range("I91").Value = 4600
$ws.Range("J91").Value = 2696
$ws.Range("K91").Value = 4600
$ws.Range("L91").Value = 2696
$ws.Range("M91").Value = -3196
$ws.Range("N91").Value = -5504

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4201
$ws.Range("I105").Value = 4201
$ws.Range("K105").Value = 4201
$ws.Range("M105").Value = -2454
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2352.5312
$ws.Range("I58").Value = 1785.1666
$ws.Range("J58").Value = 3082
$ws.Range("K58").Value = 1785.1666
$ws.Range("L58").Value = 3082
$ws.Range("M58").Value = -1582.1666
$ws.Range("N58").Value = -3488
$ws.Range("H105").Value = 944.4375
$ws.Range("I105").Value = 981.8182
$ws.Range("J105").Value = 862.2
$ws.Range("K105").Value = 981.8182
$ws.Range("L105").Value = 862.2
$ws.Range("M105").Value = 765.1818
$ws.Range("N105").Value = -4356.2
$ws.Range("H136").Value = 2352.5312
$ws.Range("I136").Value = 1785.1666
$ws.Range("J136").Value = 3082
$ws.Range("K136").Value = 5355.4998
$ws.Range("L136").Value = 9246
$ws.Range("M136").Value = -2805.4998
$ws.Range("N136").Value = -14346

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 856.5
$ws.Range("I5").Value = 757.53845
$ws.Range("K5").Value = 2272.61535
$ws.Range("M5").Value = -2160.61535
$ws.Range("H131").Value = 1032726.5
$ws.Range("J131").Value = 1236042.8
$ws.Range("L131").Value = 3708128.4
$ws.Range("N131").Value = -3718208.4
$ws.Range("H133").Value = 7622.9395
$ws.Range("J133").Value = 8138.107
$ws.Range("L133").Value = 24414.321
$ws.Range("N133").Value = -34534.321
$ws.Range("H135").Value = 856.5
$ws.Range("I135").Value = 757.53845
$ws.Range("K135").Value = 6817.84605
$ws.Range("M135").Value = -4282.84605

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2177.3235
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2177.3235
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2177.3235
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -2401.3235
$ws.Range("H21").Value = 602300
$ws.Range("I21").Value = 3000
$ws.Range("J21").Value = 2000666.6
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 2000666.6
$ws.Range("M21").Value = -2827
$ws.Range("N21").Value = -2001012.6
$ws.Range("H30").Value = 602300
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 2000666.6
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 2000666.6
$ws.Range("M30").Value = -2895
$ws.Range("N30").Value = -2000876.6
$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10302
$ws.Range("H43").Value = 22562.5
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3245.4666
$ws.Range("I81").Value = 1200
$ws.Range("J81").Value = 4609.1113
$ws.Range("K81").Value = 2400
$ws.Range("L81").Value = 9218.222599999999
$ws.Range("M81").Value = -1339
$ws.Range("N81").Value = -11340.2226
$ws.Range("H84").Value = 3245.4666
$ws.Range("I84").Value = 1200
$ws.Range("J84").Value = 4609.1113
$ws.Range("K84").Value = 12000
$ws.Range("L84").Value = 46091.113
$ws.Range("M84").Value = -6696
$ws.Range("N84").Value = -56699.113
